# Update "想去人数" (want-to-go count) figures across sheets, reflecting
# output generated at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 15119
$ws1.Range("F4").Value = 710
$ws1.Range("F5").Value = 248
$ws1.Range("F6").Value = 635
$ws1.Range("F7").Value = 1609
$ws1.Range("F8").Value = 152

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 4

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 15119
$ws4.Range("F4").Value = 710
$ws4.Range("F5").Value = 248
$ws4.Range("F8").Value = 635
$ws4.Range("F9").Value = 1609
$ws4.Range("F10").Value = 4
$ws4.Range("F11").Value = 152
